$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Swap the "Enterprises (absolute #)" row (12) and the
# "Enterprises density (per 1000 people)" row (13) so that the
# density row now appears first, followed by the absolute-count row.
# (Value2 is used because Value has a COM-bridge quirk on this runtime;
# a leading apostrophe forces numeric-looking text like "40000"/"38.7"
# to stay text instead of being coerced to a number, and resetting the
# Style back to Normal afterwards drops the transient quote-prefix
# formatting that the apostrophe trick would otherwise leave behind.)

$labelA12 = $ws.Range("A12").Value2
$valueD12 = $ws.Range("D12").Value2
$labelA13 = $ws.Range("A13").Value2
$valueD13 = $ws.Range("D13").Value2

$ws.Range("A12").Value2 = $labelA13
$ws.Range("D12").Value2 = "'" + $valueD13
$ws.Range("D12").Style = "Normal"

$ws.Range("A13").Value2 = $labelA12
$ws.Range("D13").Value2 = "'" + $valueD12
$ws.Range("D13").Style = "Normal"
